$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.04"
$ws.Range("E2").Value = "'-3.30%"
$ws.Range("D3").Value = "'40.94"
$ws.Range("E3").Value = "'-2.37%"
$ws.Range("D4").Value = "'5.045"
$ws.Range("E4").Value = "'-2.54%"
$ws.Range("D5").Value = "'0.07608"
$ws.Range("E5").Value = "'-6.00%"
$ws.Range("D6").Value = "'4.249"
$ws.Range("E6").Value = "'-2.82%"
$ws.Range("D7").Value = "'1.598"
$ws.Range("E7").Value = "'-8.88%"
$ws.Range("D8").Value = "'0.9056"
$ws.Range("E8").Value = "'-2.73%"
$ws.Range("D9").Value = "'0.09942"
$ws.Range("E9").Value = "'-11.36%"
$ws.Range("D10").Value = "'0.1769"
$ws.Range("E10").Value = "'-5.31%"
$ws.Range("D11").Value = "'0.09214"
$ws.Range("E11").Value = "'-1.20%"
$ws.Range("D12").Value = "'0.04416"
$ws.Range("E12").Value = "'-3.04%"
$ws.Range("E13").Value = "'-0.19%"
$ws.Range("D14").Value = "'0.001259"
$ws.Range("E14").Value = "'-2.97%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005801"
$ws.Range("E15").Value = "'-1.34%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.371"
$ws.Range("E16").Value = "'0.40%"
$ws.Range("B17").Value = "BTSEToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D17").Value = "'2.455"
$ws.Range("E17").Value = "'-5.13%"
$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D18").Value = "'0.3298"
$ws.Range("E18").Value = "'-1.69%"
$ws.Range("B19").Value = "MCDex"
$ws.Range("C19").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D19").Value = "'6.744"
$ws.Range("E19").Value = "'-8.76%"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").Value = "'0.1339"
$ws.Range("E20").Value = "'-2.99%"
$ws.Range("B21").Value = "ZBToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D21").Value = "'0.2843"
$ws.Range("E21").Value = "'11.49%"
$ws.Range("B22").Value = "CoinExToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D22").Value = "'0.04164"
$ws.Range("E22").Value = "'0.10%"
$ws.Range("D24").Value = "'0.004070"
$ws.Range("E24").Value = "'-6.02%"
$ws.Range("E25").Value = "'6.33%"
$ws.Range("D26").Value = "'0.0003006"
$ws.Range("E26").Value = "'0.73%"
$ws.Range("E38").Value = "'-6.34%"
$ws.Range("D39").Value = "'0.05153"
$ws.Range("E39").Value = "'-5.93%"
$ws.Range("D40").Value = "'0.007827"
$ws.Range("E40").Value = "'-2.62%"
$ws.Range("D41").Value = "'0.1306"
$ws.Range("E41").Value = "'-6.29%"
$ws.Range("D42").Value = "'0.007067"
$ws.Range("E42").Value = "'-6.56%"
$ws.Range("D43").Value = "'0.001948"
$ws.Range("E43").Value = "'-6.84%"
$ws.Range("D44").Value = "'0.008067"
$ws.Range("E44").Value = "'-2.41%"
$ws.Range("E45").Value = "'5.92%"
$ws.Range("D46").Value = "'0.00006384"
$ws.Range("E46").Value = "'-5.78%"
$ws.Range("E47").Value = "'-0.22%"
$ws.Range("E48").Value = "'-26.89%"
$ws.Range("D49").Value = "'0.006333"
$ws.Range("E49").Value = "'86.83%"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.22%"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.22%"
